$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '86.817.23'
$ws.Range('E2').Value = '  -2.91%  '
$ws.Range('D3').Value = '3.135.47'
$ws.Range('E3').Value = '  -6.69%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '203.38'
$ws.Range('E5').Value = '  -7.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '605.06'
$ws.Range('E6').Value = '  -6.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.367'
$ws.Range('E7').Value = '  -9.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.654'
$ws.Range('E8').Value = '  +6.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').Value = '3.127.10'
$ws.Range('E10').Value = '  -6.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.527'
$ws.Range('E11').Value = '  -11.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.176'
$ws.Range('E12').Value = '  +4.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000240'
$ws.Range('E13').Value = '  -17.35%  '
$ws.Range('D14').Value = '3.708.44'
$ws.Range('E14').Value = '  -7.03%  '
$ws.Range('E15').Value = '  -6.41%  '
$ws.Range('D16').Value = '86.526.74'
$ws.Range('E16').Value = '  -3.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '31.87'
$ws.Range('E17').Value = '  -13.75%  '
$ws.Range('D18').Value = '3.147.30'
$ws.Range('E18').Value = '  -6.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.95'
$ws.Range('E19').Value = '  -7.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.23'
$ws.Range('E20').Value = '  -10.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '409.24'
$ws.Range('E21').Value = '  -10.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.38'
$ws.Range('E22').Value = '  -12.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.03'
$ws.Range('E23').Value = '  -9.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.08'
$ws.Range('E24').Value = '  -9.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.65'
$ws.Range('E25').Value = '  -9.74%  '
$ws.Range('D26').Value = '3.298.38'
$ws.Range('E26').Value = '  -5.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '72.84'
$ws.Range('E27').Value = '  -7.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000128'
$ws.Range('E28').Value = '  -11.25%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.159'
$ws.Range('E30').Value = '  -22.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '533.33'
$ws.Range('E32').Value = '  -10.75%  '
$ws.Range('E33').Value = '  -12.29%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.28'
$ws.Range('E34').Value = '  -19.57%  '
$ws.Range('B35').Value = 'PancakeSwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.84'
$ws.Range('E35').Value = '  -12.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.46'
$ws.Range('E36').Value = '  -13.37%  '
$ws.Range('E37').Value = '  -8.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '21.51'
$ws.Range('E38').Value = '  -8.01%  '
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '21.77'
$ws.Range('E40').Value = '  -0.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.96'
$ws.Range('E41').Value = '  -8.34%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.88'
$ws.Range('E43').Value = '  -13.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.366'
$ws.Range('E44').Value = '  -13.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '147.75'
$ws.Range('E45').Value = '  -6.55%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '170.15'
$ws.Range('E46').Value = '  -9.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '42.61'
$ws.Range('E47').Value = '  -7.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.126'
$ws.Range('E48').Value = '  +5.06%  '
$ws.Range('E49').Value = '  -17.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.93'
$ws.Range('E50').Value = '  -13.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.579'
$ws.Range('E51').Value = '  -13.55%  '
